$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Actual_Data: append three new measurement rows (26-28)
# ---------------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("Actual_Data")

$wsData.Range("A26").Value = 96
$wsData.Range("B26").Value = 22
$wsData.Range("C26").Value = 1.23
$wsData.Range("D26").Value = "50-100"

$wsData.Range("A27").Value = 103
$wsData.Range("B27").Value = 18
$wsData.Range("C27").Value = 1.505
$wsData.Range("D27").Value = "100-150"

$wsData.Range("A28").Value = 120
$wsData.Range("B28").Value = 13
$wsData.Range("C28").Value = 1.8
$wsData.Range("D28").Value = "100-150"

# ---------------------------------------------------------------------------
# Categories: add a new "100-150" bucket in row 7, then drop the trailing
# empty placeholder rows that used to stretch down to row 12
# ---------------------------------------------------------------------------
$wsCat = $wb.Worksheets.Item("Categories")

$wsCat.Range("A7").Value = 100
$wsCat.Range("B7").Value = 150
$wsCat.Range("C7").Value = "100-150"
$wsCat.Range("D7").ClearFormats() | Out-Null
$wsCat.Range("D7").Formula = "=AVERAGE(A7,B7)"

$wsCat.Range("A8:D12").Clear() | Out-Null

# ---------------------------------------------------------------------------
# Workbook-level defined name (mirrors the MySQL-for-Excel helper name that
# Excel writes out after refreshing data from a MySQL source)
# ---------------------------------------------------------------------------
$dateFormatFormula = '=REPT(LOCAL_YEAR_FORMAT,4)&LOCAL_DATE_SEPARATOR&REPT(LOCAL_MONTH_FORMAT,2)&LOCAL_DATE_SEPARATOR&REPT(LOCAL_DAY_FORMAT,2)&" "&REPT(LOCAL_HOUR_FORMAT,2)&LOCAL_TIME_SEPARATOR&REPT(LOCAL_MINUTE_FORMAT,2)&LOCAL_TIME_SEPARATOR&REPT(LOCAL_SECOND_FORMAT,2)'
$definedName = $wb.Names.Add("LOCAL_MYSQL_DATE_FORMAT", $dateFormatFormula)
$definedName.Visible = $false

# ---------------------------------------------------------------------------
# Selections: leave Categories' active cell on D8, and come back to rest on
# Actual_Data (which stays the selected tab) at H27
# ---------------------------------------------------------------------------
$wsCat.Range("D8").Select() | Out-Null
$wsData.Range("H27").Select() | Out-Null

Write-Output "done"
